$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Agosto de 2020 a las 19:45"

# Country label reshuffle: Guatemala overtakes Bielorrusia in ranking (rows 43-44)
$ws.Range("A43").Value = "Guatemala"
$ws.Range("A44").Value = "Bielorrusia"

# Country label reshuffle: Marruecos overtakes Portugal/Honduras/Singapur (rows 49-52)
$ws.Range("A49").Value = "Marruecos"
$ws.Range("A50").Value = "Portugal"
$ws.Range("A51").Value = "Honduras"
$ws.Range("A52").Value = "Singapur"

# Updated case statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 6016241
$ws.Range("C4").Value = 15876
$ws.Range("D4").Value = 3320336
$ws.Range("E4").Value = 2511933
$ws.Range("G4").Value = 319
$ws.Range("H4").Value = 183972
$ws.Range("B6").Value = 3382152
$ws.Range("C6").Value = 74403
$ws.Range("D6").Value = 2582179
$ws.Range("E6").Value = 738298
$ws.Range("G6").Value = 1046
$ws.Range("H6").Value = 61675
$ws.Range("B12").Value = 430599
$ws.Range("C12").Value = 3781
$ws.Range("G12").Value = 25
$ws.Range("H12").Value = 28996
$ws.Range("B23").Value = 240030
$ws.Range("C23").Value = 1030
$ws.Range("E23").Value = 16438
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 9359
$ws.Range("B24").Value = 219435
$ws.Range("C24").Value = 3651
$ws.Range("D24").Value = 161009
$ws.Range("E24").Value = 51686
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = 6740
$ws.Range("B27").Value = 126646
$ws.Range("C27").Value = 229
$ws.Range("D27").Value = 112619
$ws.Range("E27").Value = 4929
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 9098
$ws.Range("B31").Value = 111219
$ws.Range("C31").Value = 670
$ws.Range("D31").Value = 95202
$ws.Range("E31").Value = 9546
$ws.Range("G31").Value = 61
$ws.Range("H31").Value = 6471
$ws.Range("B43").Value = 71856
$ws.Range("C43").Value = 1142
$ws.Range("D43").Value = 59641
$ws.Range("E43").Value = 9530
$ws.Range("G43").Value = 23
$ws.Range("H43").Value = 2685
$ws.Range("B44").Value = 71165
$ws.Range("C44").Value = 191
$ws.Range("D44").Value = 69650
$ws.Range("E44").Value = 853
$ws.Range("G44").Value = 5
$ws.Range("H44").Value = 662
$ws.Range("B49").Value = 57085
$ws.Range("C49").Value = 1221
$ws.Range("D49").Value = 41901
$ws.Range("E49").Value = 14173
$ws.Range("G49").Value = 27
$ws.Range("H49").Value = 1011
$ws.Range("B50").Value = 56673
$ws.Range("C50").Value = 399
$ws.Range("D50").Value = 41357
$ws.Range("E50").Value = 13507
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 1809
$ws.Range("B51").Value = 56649
$ws.Range("C51").Value = 772
$ws.Range("D51").Value = 9169
$ws.Range("E51").Value = 45733
$ws.Range("G51").Value = 44
$ws.Range("H51").Value = 1747
$ws.Range("B52").Value = 56572
$ws.Range("C52").Value = 77
$ws.Range("D52").Value = 55139
$ws.Range("E52").Value = 1406
$ws.Range("H52").Value = 27
$ws.Range("E54").Value = 3157
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 187
$ws.Range("B55").Value = 46407
$ws.Range("C55").Value = 1186
$ws.Range("D55").Value = 16829
$ws.Range("E55").Value = 28833
$ws.Range("G55").Value = 20
$ws.Range("H55").Value = 745
$ws.Range("B56").Value = 43841
$ws.Range("C56").Value = 72
$ws.Range("D56").Value = 42246
$ws.Range("E56").Value = 1325
$ws.Range("B59").Value = 43016
$ws.Range("C59").Value = 397
$ws.Range("D59").Value = 30157
$ws.Range("E59").Value = 11384
$ws.Range("G59").Value = 10
$ws.Range("H59").Value = 1475
$ws.Range("B70").Value = 28453
$ws.Range("C70").Value = 90
$ws.Range("E70").Value = 3312
$ws.Range("B74").Value = 23169
$ws.Range("C74").Value = 218
$ws.Range("D74").Value = 17190
$ws.Range("E74").Value = 5561
$ws.Range("B132").Value = 2679
$ws.Range("C132").Value = 155
$ws.Range("E132").Value = 1383
$ws.Range("G132").Value = 2
$ws.Range("H132").Value = 28
$ws.Range("B141").Value = 1933
$ws.Range("C141").Value = 3
$ws.Range("D141").Value = 1098
$ws.Range("G141").Value = 2
$ws.Range("H141").Value = 562
